# Applies the "ExcelDatabase" data-refresh edit:
#   - renames several dummy customer companies
#   - adds a new customer ("Dummy"), worker ("John Doe"), two services
#     ("Test Service", "Boris"), a new project (row 7) and new bill rows
#   - updates related bridge tables (worker<->project, project<->services,
#     project<->bill) and the bill sheet totals / addresses

$wb = $excel.ActiveWorkbook

function Set-Text {
    param($range, [string]$text)
    # Force the cell to be stored as a text/shared-string value even when the
    # text looks like a number or a date (e.g. phone numbers, ISO dates).
    $range.Formula = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# customer sheet
# ---------------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item("customer")

$wsCustomer.Cells.Item(2, 2).Value = "One Company"
$wsCustomer.Cells.Item(3, 2).Value = "Company Two"
$wsCustomer.Cells.Item(4, 2).Value = "Third Company"
$wsCustomer.Cells.Item(5, 2).Value = "4Company"
$wsCustomer.Cells.Item(6, 2).Value = "FiveGuys Burger"
$wsCustomer.Cells.Item(7, 2).Value = "6th and Co."

# new customer row
$wsCustomer.Cells.Item(8, 1).Value = 7
$wsCustomer.Cells.Item(8, 2).Value = "Dummy"
Set-Text $wsCustomer.Cells.Item(8, 3) "44754895625"
$wsCustomer.Cells.Item(8, 4).Value = "dummy@testcorp.com"
$wsCustomer.Cells.Item(8, 5).Value = 119
$wsCustomer.Cells.Item(8, 6).Value = 1.5

# ---------------------------------------------------------------------
# project sheet
# ---------------------------------------------------------------------
$wsProject = $wb.Worksheets.Item("project")

$wsProject.Cells.Item(2, 3).Value = 9
$wsProject.Cells.Item(2, 23).Value = 50

# new project row
$wsProject.Cells.Item(8, 1).Value = 7
Set-Text $wsProject.Cells.Item(8, 2) "2022-02-16"
$wsProject.Cells.Item(8, 3).Value = "This is some test data filling in the specimen procedure information. Please hold till a member of our team is able to speak with you..."
for ($col = 4; $col -le 10; $col++) {
    $wsProject.Cells.Item(8, $col).Value = "Test Data...."
}
$wsProject.Cells.Item(8, 11).Value = "Yes"
for ($col = 12; $col -le 21; $col++) {
    $wsProject.Cells.Item(8, $col).Value = "Test Data...."
}
$wsProject.Cells.Item(8, 22).Value = 2
$wsProject.Cells.Item(8, 23).Value = 4050
$wsProject.Cells.Item(8, 24).Value = 7

# ---------------------------------------------------------------------
# services sheet
# ---------------------------------------------------------------------
$wsServices = $wb.Worksheets.Item("services")

$wsServices.Cells.Item(15, 1).Value = 15
$wsServices.Cells.Item(15, 2).Value = "Test Service"
$wsServices.Cells.Item(15, 3).Value = 50
$wsServices.Cells.Item(15, 4).Value = 25
$wsServices.Cells.Item(15, 5).Value = 75
$wsServices.Cells.Item(15, 6).Value = "hour"

$wsServices.Cells.Item(16, 1).Value = 16
$wsServices.Cells.Item(16, 2).Value = "Boris"
$wsServices.Cells.Item(16, 3).Value = 300
$wsServices.Cells.Item(16, 4).Value = 150
$wsServices.Cells.Item(16, 5).Value = 450
$wsServices.Cells.Item(16, 6).Value = "sample"

# ---------------------------------------------------------------------
# worker sheet
# ---------------------------------------------------------------------
$wsWorker = $wb.Worksheets.Item("worker")

$wsWorker.Cells.Item(8, 1).Value = 7
$wsWorker.Cells.Item(8, 2).Value = "John Doe"
Set-Text $wsWorker.Cells.Item(8, 3) "44765496216"
$wsWorker.Cells.Item(8, 4).Value = "john@doe.com"
$wsWorker.Cells.Item(8, 5).Value = 7

# ---------------------------------------------------------------------
# workerprojectbridge sheet
# ---------------------------------------------------------------------
$wsWpb = $wb.Worksheets.Item("workerprojectbridge")

$wsWpb.Cells.Item(7, 1).Value = 6
$wsWpb.Cells.Item(7, 2).Value = 6
$wsWpb.Cells.Item(7, 3).Value = 6

$wsWpb.Cells.Item(8, 1).Value = 7
$wsWpb.Cells.Item(8, 2).Value = 7
$wsWpb.Cells.Item(8, 3).Value = 7

# ---------------------------------------------------------------------
# projectservicesbridge sheet
# ---------------------------------------------------------------------
$wsPsb = $wb.Worksheets.Item("projectservicesbridge")

$wsPsb.Cells.Item(8, 2).Value = 5
$wsPsb.Cells.Item(8, 3).Value = 50
$wsPsb.Cells.Item(8, 5).Value = 6

$wsPsb.Cells.Item(9, 1).Value = 8
$wsPsb.Cells.Item(9, 2).Value = 5
$wsPsb.Cells.Item(9, 3).Value = 200
$wsPsb.Cells.Item(9, 4).Value = 1
$wsPsb.Cells.Item(9, 5).Value = 4

$wsPsb.Cells.Item(10, 1).Value = 9
$wsPsb.Cells.Item(10, 2).Value = 5
$wsPsb.Cells.Item(10, 3).Value = 200
$wsPsb.Cells.Item(10, 4).Value = 1
$wsPsb.Cells.Item(10, 5).Value = 10

$wsPsb.Cells.Item(11, 1).Value = 10
$wsPsb.Cells.Item(11, 2).Value = 2
$wsPsb.Cells.Item(11, 3).Value = 160
$wsPsb.Cells.Item(11, 4).Value = 1
$wsPsb.Cells.Item(11, 5).Value = 14

$wsPsb.Cells.Item(12, 1).Value = 11
$wsPsb.Cells.Item(12, 2).Value = 5
$wsPsb.Cells.Item(12, 3).Value = 50
$wsPsb.Cells.Item(12, 4).Value = 1
$wsPsb.Cells.Item(12, 5).Value = 6

$wsPsb.Cells.Item(13, 1).Value = 12
$wsPsb.Cells.Item(13, 2).Value = 9
$wsPsb.Cells.Item(13, 3).Value = 4050
$wsPsb.Cells.Item(13, 4).Value = 7
$wsPsb.Cells.Item(13, 5).Value = 16

# ---------------------------------------------------------------------
# projectbillbridge sheet
# ---------------------------------------------------------------------
$wsPbb = $wb.Worksheets.Item("projectbillbridge")

$wsPbb.Cells.Item(8, 1).Value = 7
$wsPbb.Cells.Item(8, 2).Value = 3
$wsPbb.Cells.Item(8, 3).Value = 4

$wsPbb.Cells.Item(9, 1).Value = 8
$wsPbb.Cells.Item(9, 2).Value = 4
$wsPbb.Cells.Item(9, 3).Value = 7

# ---------------------------------------------------------------------
# bill sheet
# ---------------------------------------------------------------------
$wsBill = $wb.Worksheets.Item("bill")

# row 2 (bill_id 1): new address, drop the extra1 name/cost, new total
$wsBill.Cells.Item(2, 3).Value = "143, Fake street, Glasg"
$wsBill.Cells.Item(2, 4).ClearContents()
$wsBill.Cells.Item(2, 5).ClearContents()
$wsBill.Cells.Item(2, 8).Value = 180

# row 3 (bill_id 2): drop the extra1/extra2 name/cost, new total
$wsBill.Cells.Item(3, 4).ClearContents()
$wsBill.Cells.Item(3, 5).ClearContents()
$wsBill.Cells.Item(3, 6).ClearContents()
$wsBill.Cells.Item(3, 7).ClearContents()
$wsBill.Cells.Item(3, 8).Value = 292.5

# row 4 (bill_id 3): new bill
$wsBill.Cells.Item(4, 1).Value = 3
Set-Text $wsBill.Cells.Item(4, 2) "2022-02-15"
$wsBill.Cells.Item(4, 3).Value = "Address, GL12 3BC"
$wsBill.Cells.Item(4, 8).Value = 82.5
$wsBill.Cells.Item(4, 9).Value = 5

# row 5 (bill_id 4): new bill with extra services
$wsBill.Cells.Item(5, 1).Value = 4
Set-Text $wsBill.Cells.Item(5, 2) "2022-02-16"
$wsBill.Cells.Item(5, 3).Value = "Address, GL12 3BC"
$wsBill.Cells.Item(5, 4).Value = "First Service"
$wsBill.Cells.Item(5, 5).Value = 10
$wsBill.Cells.Item(5, 6).Value = "Second Service"
$wsBill.Cells.Item(5, 7).Value = 25
$wsBill.Cells.Item(5, 8).Value = 4085
$wsBill.Cells.Item(5, 9).Value = 7
